# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.251.53'
$ws.Range('E2').Value = '  -0.72%  '
$ws.Range('D3').Value = '3.309.65'
$ws.Range('E3').Value = '  -2.14%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '189.66'
$ws.Range('E5').Value = '  +2.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '562.20'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('E8').Value = '  -2.06%  '
$ws.Range('D9').Value = '3.301.32'
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.185'
$ws.Range('E10').Value = '  -1.91%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.588'
$ws.Range('E11').Value = '  -1.75%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '47.89'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('E13').Value = '  +0.54%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.73'
$ws.Range('E14').Value = '  -0.68%  '
$ws.Range('D15').Value = '3.840.23'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '614.83'
$ws.Range('E16').Value = '  +0.96%  '
$ws.Range('B17').Value = 'Chainlink'
$ws.Range('C17').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '18.10'
$ws.Range('E17').Value = '  -1.53%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Value = '66.287.63'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').Value = '3.318.64'
$ws.Range('E20').Value = '  -2.09%  '
$ws.Range('E21').Value = '  -4.44%  '
$ws.Range('E22').Value = '  -1.05%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.50'
$ws.Range('E23').Value = '  +8.57%  '
$ws.Range('E24').Value = '  -2.06%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '101.77'
$ws.Range('E26').Value = '  -2.24%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.01'
$ws.Range('E27').Value = '  -0.08%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.76'
$ws.Range('E28').Value = '  +0.79%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.79'
$ws.Range('E29').Value = '  +2.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '8.63'
$ws.Range('E30').Value = '  -2.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '30.36'
$ws.Range('E31').Value = '  -1.82%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.77'
$ws.Range('E32').Value = '  +6.22%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.08'
$ws.Range('E33').Value = '  +4.53%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '577.34'
$ws.Range('E34').Value = '  +3.34%  '
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('D37').Value = '3.739.65'
$ws.Range('E37').Value = '  -3.99%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '57.31'
$ws.Range('E38').Value = '  -1.64%  '
$ws.Range('E39').Value = '  -0.11%  '
$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0736'
$ws.Range('E40').Value = '  +1.19%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.35'
$ws.Range('E41').Value = '  -3.78%  '
$ws.Range('B42').Value = 'CoreDAO'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.50'
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '34.15'
$ws.Range('E43').Value = '  +4.72%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.130'
$ws.Range('E44').Value = '  +0.84%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.74'
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.342'
$ws.Range('E46').Value = '  -2.98%  '
$ws.Range('E47').Value = '  +1.41%  '
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.129'
$ws.Range('E49').Value = '  -1.52%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.59'
$ws.Range('E50').Value = '  -4.07%  '
$ws.Range('E51').Value = '  +0.03%  '
